$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 17.01.2022 01:15"

# Row 5 (Makro): swap price values and change Delta/Date cells to text
$ws.Range("B5").Value = 34.5
$ws.Range("C5").Value = 33.9

# Force D5/E5 to be stored as text (not numbers/dates), then restore the
# default (unformatted) style so no lingering number format is left behind.
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "+0.6"
$ws.Range("E5").Value = "2022-01-17 01:15:08"
$ws.Range("D5:E5").Style = "Normal"
